$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task Sprint 2")

# Aggiornamento Stato dei task assegnati (Sprint 2): da "Non Iniziata" a "Completata"
$ws.Range("E18").Value = "Completata"
$ws.Range("E19").Value = "Completata"
$ws.Range("E21").Value = "Completata"
$ws.Range("E23").Value = "Completata"

# Allinea la formattazione delle celle di stato/colonna D aggiornate allo
# stile già usato per le righe "Completata" (es. E18), cosi' da riutilizzare
# lo stesso formato invece di crearne uno nuovo.
$ws.Range("E18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("E23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
